$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "87.943.64"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "3.112.10"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'213.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").Value = "'634.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.45%  "
$ws.Range("D7").Value = "'0.387"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").Value = "'0.833"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +21.38%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "3.108.87"
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("D11").Value = "'0.575"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "'5.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.95%  "
$ws.Range("D15").Value = "87.950.15"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").Value = "3.684.84"
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("D17").Value = "'31.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("D18").Value = "3.110.82"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").Value = "'3.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.10%  "
$ws.Range("D20").Value = "'0.0000221"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +18.62%  "
$ws.Range("D21").Value = "'13.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").Value = "'423.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("D23").Value = "'8.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("D24").Value = "'4.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.43%  "
$ws.Range("D25").Value = "'5.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.87%  "
$ws.Range("D26").Value = "'82.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.32%  "
$ws.Range("D27").Value = "'11.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").Value = "3.275.42"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  -6.64%  "
$ws.Range("D32").Value = "'3.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("D33").Value = "'8.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("E34").Value = "  +17.78%  "
$ws.Range("D35").Value = "'499.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.46%  "
$ws.Range("D36").Value = "'6.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.03%  "
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").Value = "'22.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.66%  "
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.138"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.08%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.95%  "
$ws.Range("D46").Value = "'146.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.21%  "
$ws.Range("D47").Value = "'43.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0655"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +12.06%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'162.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.26%  "
$ws.Range("D50").Value = "'0.716"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("E51").Value = "  -3.06%  "
